# Update the 25 two-digit multiplication expressions in the table
# (5 data rows x 5 columns) to the new values from the commit.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
if ($cell.Range.Text -like "18×22=*") {
    $cell.Range.Text = "51×55="
}

$cell = $t.Cell(1, 2)
if ($cell.Range.Text -like "95×53=*") {
    $cell.Range.Text = "78×45="
}

$cell = $t.Cell(1, 3)
if ($cell.Range.Text -like "55×71=*") {
    $cell.Range.Text = "27×33="
}

$cell = $t.Cell(1, 4)
if ($cell.Range.Text -like "43×51=*") {
    $cell.Range.Text = "88×90="
}

$cell = $t.Cell(1, 5)
if ($cell.Range.Text -like "56×74=*") {
    $cell.Range.Text = "95×97="
}

$cell = $t.Cell(5, 1)
if ($cell.Range.Text -like "62×66=*") {
    $cell.Range.Text = "28×45="
}

$cell = $t.Cell(5, 2)
if ($cell.Range.Text -like "36×25=*") {
    $cell.Range.Text = "30×65="
}

$cell = $t.Cell(5, 3)
if ($cell.Range.Text -like "54×25=*") {
    $cell.Range.Text = "26×88="
}

$cell = $t.Cell(5, 4)
if ($cell.Range.Text -like "90×52=*") {
    $cell.Range.Text = "90×55="
}

$cell = $t.Cell(5, 5)
if ($cell.Range.Text -like "95×53=*") {
    $cell.Range.Text = "38×22="
}

$cell = $t.Cell(10, 1)
if ($cell.Range.Text -like "82×92=*") {
    $cell.Range.Text = "46×69="
}

$cell = $t.Cell(10, 2)
if ($cell.Range.Text -like "90×99=*") {
    $cell.Range.Text = "68×25="
}

$cell = $t.Cell(10, 3)
if ($cell.Range.Text -like "27×78=*") {
    $cell.Range.Text = "22×61="
}

$cell = $t.Cell(10, 4)
if ($cell.Range.Text -like "53×86=*") {
    $cell.Range.Text = "73×38="
}

$cell = $t.Cell(10, 5)
if ($cell.Range.Text -like "12×86=*") {
    $cell.Range.Text = "65×24="
}

$cell = $t.Cell(15, 1)
if ($cell.Range.Text -like "40×66=*") {
    $cell.Range.Text = "51×13="
}

$cell = $t.Cell(15, 2)
if ($cell.Range.Text -like "22×30=*") {
    $cell.Range.Text = "54×19="
}

$cell = $t.Cell(15, 3)
if ($cell.Range.Text -like "99×86=*") {
    $cell.Range.Text = "77×89="
}

$cell = $t.Cell(15, 4)
if ($cell.Range.Text -like "25×33=*") {
    $cell.Range.Text = "12×72="
}

$cell = $t.Cell(15, 5)
if ($cell.Range.Text -like "66×33=*") {
    $cell.Range.Text = "33×31="
}

$cell = $t.Cell(20, 1)
if ($cell.Range.Text -like "96×18=*") {
    $cell.Range.Text = "78×90="
}

$cell = $t.Cell(20, 2)
if ($cell.Range.Text -like "54×92=*") {
    $cell.Range.Text = "92×73="
}

$cell = $t.Cell(20, 3)
if ($cell.Range.Text -like "85×95=*") {
    $cell.Range.Text = "93×67="
}

$cell = $t.Cell(20, 4)
if ($cell.Range.Text -like "53×66=*") {
    $cell.Range.Text = "59×85="
}

$cell = $t.Cell(20, 5)
if ($cell.Range.Text -like "21×75=*") {
    $cell.Range.Text = "75×38="
}

